$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) holds numeric-looking strings (e.g. "309.60", "0.0300") that
# Excel would silently coerce to numbers (dropping trailing zeros / dot-grouping).
# Force the whole Price column to Text format first so assigned strings are kept verbatim.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "42.650.80"
$ws.Range("E2").Value = "  -0.37%  "
$ws.Range("D3").Value = "2.532.44"
$ws.Range("E3").Value = "  -1.27%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "309.60"
$ws.Range("E5").Value = "  -1.30%  "
$ws.Range("D6").Value = "100.22"
$ws.Range("E6").Value = "  +0.88%  "
$ws.Range("D7").Value = "0.568"
$ws.Range("E7").Value = "  -1.23%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "0.523"
$ws.Range("E9").Value = "  -2.03%  "
$ws.Range("D10").Value = "35.78"
$ws.Range("E10").Value = "  -0.21%  "
$ws.Range("E11").Value = "  -0.80%  "
$ws.Range("E12").Value = "  -1.63%  "
$ws.Range("E13").Value = "  +0.78%  "
$ws.Range("D14").Value = "2.921.44"
$ws.Range("E14").Value = "  -1.40%  "
$ws.Range("E15").Value = "  -2.96%  "
$ws.Range("D16").Value = "2.520.74"
$ws.Range("E16").Value = "  -2.66%  "
$ws.Range("E17").Value = "  -4.17%  "
$ws.Range("D18").Value = "42.643.69"
$ws.Range("E18").Value = "  -0.50%  "
$ws.Range("E19").Value = "  -0.97%  "
$ws.Range("D21").Value = "12.25"
$ws.Range("E21").Value = "  -2.43%  "
$ws.Range("D22").Value = "69.28"
$ws.Range("E22").Value = "  -0.16%  "
$ws.Range("D23").Value = "243.38"
$ws.Range("E23").Value = "  -2.30%  "
$ws.Range("D24").Value = "2.87"
$ws.Range("E24").Value = "  -2.78%  "
$ws.Range("E25").Value = "  -1.84%  "
$ws.Range("E26").Value = "  +0.71%  "
$ws.Range("D27").Value = "25.49"
$ws.Range("E27").Value = "  -5.37%  "
$ws.Range("E28").Value = "  -1.96%  "
$ws.Range("E29").Value = "  -0.88%  "
$ws.Range("D30").Value = "38.53"
$ws.Range("E30").Value = "  -4.48%  "
$ws.Range("D31").Value = "157.97"
$ws.Range("E31").Value = "  +0.47%  "
$ws.Range("E33").Value = "  +9.84%  "
$ws.Range("E34").Value = "  -1.48%  "
$ws.Range("D35").Value = "0.0784"
$ws.Range("E35").Value = "  -2.03%  "
$ws.Range("D36").Value = "18.29"
$ws.Range("E36").Value = "  -2.50%  "
$ws.Range("E37").Value = "  -7.13%  "
$ws.Range("D38").Value = "1.97"
$ws.Range("E38").Value = "  -6.85%  "
$ws.Range("D39").Value = "0.111"
$ws.Range("E39").Value = "  -1.22%  "
$ws.Range("E40").Value = "  -0.87%  "
$ws.Range("D41").Value = "4.27"
$ws.Range("E41").Value = "  +3.99%  "
$ws.Range("D42").Value = "22.54"
$ws.Range("E42").Value = "  -4.16%  "
$ws.Range("E43").Value = "  +0.15%  "
$ws.Range("D44").Value = "0.0300"
$ws.Range("E44").Value = "  -0.85%  "
$ws.Range("E45").Value = "  +1.54%  "
$ws.Range("D46").Value = "1.992.03"
$ws.Range("E46").Value = "  -0.69%  "
$ws.Range("E47").Value = "  -0.46%  "
$ws.Range("D48").Value = "2.774.57"
$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").Value = "0.190"
$ws.Range("E49").Value = "  -3.20%  "
$ws.Range("B50").Value = "BitcoinSV"
$ws.Range("C50").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D50").Value = "79.47"
$ws.Range("E50").Value = "  -2.67%  "
$ws.Range("D51").Value = "72.16"
$ws.Range("E51").Value = "  -3.31%  "
